$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Michael Fay Park, Longford Slashers", "Longford"),
    @("Moneygall", "Tipperary"),
    @("Azzurri Walsh Park, Waterford", "Waterford"),
    @("Rurai Og Cushendall", "Antrim")
)

$startRow = 121
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
